$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.508.46'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.252.37'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.30%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '233.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.634'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.91%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '69.90'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.563'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.10%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0997'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '58.61'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '36.25'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +12.65%  '
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.76'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.587.50'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.29%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.16'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.22%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.861'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.244.18'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '42.345.60'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.87%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0980'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.88%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.56'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.26%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.82'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.07%  '
$ws.Range('E24').Value = '  +5.84%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.42'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('E28').Value = '  -1.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.21'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.47'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.35%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.61'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -6.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.121'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.68%  '
$ws.Range('E33').Value = '  -4.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0731'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.42'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.79%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.72'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.65'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '21.67'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +16.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.28'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.02'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0271'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '65.87'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.97'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -11.48%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.12'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.193'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.53'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +13.43%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.18'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.11'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +9.74%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.34'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.34%  '
